$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace the old 3-column "Name/Position/Skill" table with the new
#     8-column Italian "calcetto" player-ratings table -------------------
$ws.Range("A1:H14").ClearContents()

# Header row — written in this particular column order so the shared-string
# table comes out in the same sequence as the authored workbook.
$ws.Range("B1").Value = "Difesa"
$ws.Range("C1").Value = "Attacco"
$ws.Range("G1").Value = "Porta"
$ws.Range("E1").Value = "Velocità"
$ws.Range("F1").Value = "Tiro"
$ws.Range("D1").Value = "Centrocampo"
$ws.Range("A1").Value = "Nome"

# Player names (column A), in row order except Antonello (row 8) which was
# appended to the shared-string table after the rest of the roster.
$ws.Range("A2").Value = "Pietro"
$ws.Range("A3").Value = "Simone"
$ws.Range("A4").Value = "Luciano"
$ws.Range("A5").Value = "Marco Ma"
$ws.Range("A6").Value = "Marco Me"
$ws.Range("A7").Value = "Oscar"
$ws.Range("A9").Value = "Roberto M"
$ws.Range("A10").Value = "Roberto"
$ws.Range("A11").Value = "Max"
$ws.Range("A8").Value = "Antonello"
$ws.Range("A12").Value = "Fabrizio"
$ws.Range("A13").Value = "Gianfranco"
$ws.Range("A14").Value = "Andrea "

$ws.Range("H1").Value = "Presenza"

# Numeric rating grid (B2:G14) + Presenza column (H2:H14).
$grid = New-Object 'object[,]' 13,7
$grid[0,0]=3; $grid[0,1]=3; $grid[0,2]=3; $grid[0,3]=3; $grid[0,4]=2; $grid[0,5]=0; $grid[0,6]=1
$grid[1,0]=2; $grid[1,1]=3; $grid[1,2]=3; $grid[1,3]=2; $grid[1,4]=3; $grid[1,5]=1; $grid[1,6]=1
$grid[2,0]=2; $grid[2,1]=3; $grid[2,2]=3; $grid[2,3]=3; $grid[2,4]=2; $grid[2,5]=0; $grid[2,6]=1
$grid[3,0]=3; $grid[3,1]=1; $grid[3,2]=1; $grid[3,3]=2; $grid[3,4]=1; $grid[3,5]=0; $grid[3,6]=1
$grid[4,0]=3; $grid[4,1]=2; $grid[4,2]=3; $grid[4,3]=2; $grid[4,4]=2; $grid[4,5]=0; $grid[4,6]=1
$grid[5,0]=3; $grid[5,1]=2; $grid[5,2]=1; $grid[5,3]=2; $grid[5,4]=1; $grid[5,5]=0; $grid[5,6]=0
$grid[6,0]=2; $grid[6,1]=1; $grid[6,2]=1; $grid[6,3]=1; $grid[6,4]=1; $grid[6,5]=1; $grid[6,6]=1
$grid[7,0]=3; $grid[7,1]=3; $grid[7,2]=3; $grid[7,3]=2; $grid[7,4]=3; $grid[7,5]=0; $grid[7,6]=0
$grid[8,0]=2; $grid[8,1]=2; $grid[8,2]=2; $grid[8,3]=2; $grid[8,4]=2; $grid[8,5]=0; $grid[8,6]=1
$grid[9,0]=1; $grid[9,1]=3; $grid[9,2]=2; $grid[9,3]=2; $grid[9,4]=2; $grid[9,5]=0; $grid[9,6]=1
$grid[10,0]=3; $grid[10,1]=1; $grid[10,2]=2; $grid[10,3]=1; $grid[10,4]=2; $grid[10,5]=0; $grid[10,6]=1
$grid[11,0]=3; $grid[11,1]=1; $grid[11,2]=3; $grid[11,3]=2; $grid[11,4]=2; $grid[11,5]=0; $grid[11,6]=1
$grid[12,0]=2; $grid[12,1]=1; $grid[12,2]=1; $grid[12,3]=1; $grid[12,4]=2; $grid[12,5]=0; $grid[12,6]=0

$ws.Range("B2:H14").Value = $grid

# --- Formatting --------------------------------------------------------
# A1 ("Nome") is bold, not centered.
$ws.Range("A1").Font.Bold = $true

# B1:H1 are bold + centered.
$ws.Range("B1:H1").Font.Bold = $true
$ws.Range("B1:H1").HorizontalAlignment = -4108   # xlCenter

# All the numeric rating cells (B2:H14) are centered.
$ws.Range("B2:H14").HorizontalAlignment = -4108  # xlCenter

# --- Column widths -------------------------------------------------------
# Column D ("Centrocampo") needs extra width to fit the header text.
$ws.Columns.Item(4).AutoFit()

# --- Selection -----------------------------------------------------------
[void]$ws.Range("H13").Select()
